$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.150.99"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.926.82"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Formula = "'591.84"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Formula = "'145.18"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Formula = "'7.00"
$ws.Range("E9").Value = "  +5.38%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Formula = "'33.79"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "3.412.81"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "61.059.77"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "2.929.56"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Formula = "'435.69"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Formula = "'13.45"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Formula = "'81.44"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").Formula = "'2.21"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Formula = "'11.89"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Formula = "'2.60"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("D32").Formula = "'26.73"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "0.0₃0867"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Formula = "'42.04"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Formula = "'376.02"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "2.687.91"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Formula = "'133.59"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Formula = "'24.01"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("E51").Value = "  +0.00%  "
